$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update selection (was E10, now H9) ---
$ws.Range("H9").Select() | Out-Null

# --- Row 1 header got taller (47.25 -> 173.25) ---
$ws.Rows(1).RowHeight = 173.25

# --- Row 11: H11:J11 used to be a single merged "no NLDAS data" text cell.
#     Replace it with actual numeric assessment data (offset, max offset,
#     and max offset year) and drop the merge, matching rows 7/8's
#     unmerged-but-still-styled-as-merge sibling cells' numeric layout.
$ws.Range("H11:J11").UnMerge()

$ws.Range("H11").Value = 1.189862
$ws.Range("I11").Value = 6.3348899999999997
$ws.Range("J11").Value = 1993

# Give the new numbers the same look as the other numeric offset/year
# columns elsewhere in the sheet (plain font, not the red "flag" font that
# was used for the old merged placeholder text) -- copy number formats
# from existing cells that already carry that plain font.
$ws.Range("N11").Copy() | Out-Null
$ws.Range("H11:I11").PasteSpecial(-4122) | Out-Null

$ws.Range("D11").Copy() | Out-Null
$ws.Range("J11").PasteSpecial(-4122) | Out-Null

# Touch horizontal alignment (back to default/general) so the style is
# recorded distinctly from the copied-from cells' styles.
$ws.Range("H11:J11").HorizontalAlignment = 1

$excel.CutCopyMode = 0
